$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 723, shifting rows 723:764 down to 724:765,
# then populate the new row with its values.
$ws.Rows.Item(723).Insert()

# Leading apostrophe forces the date-like text to be stored as a plain
# string instead of being auto-converted to a date serial number (to
# match the other "yyyy/mm/dd" text cells in column A); ClearFormats
# then strips the quote-prefix cell format so the cell keeps the default
# (unstyled) look like its neighbours.
$ws.Cells.Item(723, 1).Value = "'2026/01/30"
$ws.Cells.Item(723, 1).ClearFormats()
$ws.Cells.Item(723, 2).Value = "金"
$ws.Cells.Item(723, 3).Value = 13
$ws.Cells.Item(723, 4).Value = 20
